$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 1661.3529
$ws.Cells.Item(29, 9).Value = 250
$ws.Cells.Item(29, 11).Value = 750
$ws.Cells.Item(29, 13).Value = -469
$ws.Cells.Item(33, 8).Value = 239.66667
$ws.Cells.Item(33, 9).Value = 133.54546
$ws.Cells.Item(33, 10).Value = 531.5
$ws.Cells.Item(33, 11).Value = 133.54546
$ws.Cells.Item(33, 12).Value = 531.5
$ws.Cells.Item(33, 13).Value = 95.45454000000001
$ws.Cells.Item(33, 14).Value = -989.5
$ws.Cells.Item(38, 8).Value = 643.7143
$ws.Cells.Item(38, 9).Value = 301.4
$ws.Cells.Item(38, 10).Value = 1499.5
$ws.Cells.Item(38, 11).Value = 904.1999999999999
$ws.Cells.Item(38, 12).Value = 4498.5
$ws.Cells.Item(38, 13).Value = -532.1999999999999
$ws.Cells.Item(38, 14).Value = -5242.5
$ws.Cells.Item(64, 8).Value = 4494.4443
$ws.Cells.Item(64, 9).Value = 4747.5
$ws.Cells.Item(64, 11).Value = 4747.5
$ws.Cells.Item(64, 13).Value = -4499.5
$ws.Cells.Item(67, 8).Value = 4494.4443
$ws.Cells.Item(67, 9).Value = 4747.5
$ws.Cells.Item(67, 11).Value = 4747.5
$ws.Cells.Item(67, 13).Value = -3889.5
$ws.Cells.Item(70, 8).Value = 2073.6667
$ws.Cells.Item(70, 9).Value = 2365
$ws.Cells.Item(70, 10).Value = 1840.6
$ws.Cells.Item(70, 11).Value = 7095
$ws.Cells.Item(70, 12).Value = 5521.799999999999
$ws.Cells.Item(70, 13).Value = -6825
$ws.Cells.Item(70, 14).Value = -6061.799999999999
$ws.Cells.Item(73, 8).Value = 2073.6667
$ws.Cells.Item(73, 9).Value = 2365
$ws.Cells.Item(73, 10).Value = 1840.6
$ws.Cells.Item(73, 11).Value = 7095
$ws.Cells.Item(73, 12).Value = 5521.799999999999
$ws.Cells.Item(73, 13).Value = -6159
$ws.Cells.Item(73, 14).Value = -7393.799999999999
$ws.Cells.Item(137, 8).Value = 1042.4103
$ws.Cells.Item(137, 9).Value = 815.451
$ws.Cells.Item(137, 11).Value = 2446.353
$ws.Cells.Item(137, 13).Value = 103.6469999999999
$ws.Cells.Item(141, 8).Value = 664.8461
$ws.Cells.Item(141, 9).Value = 571.4400000000001
$ws.Cells.Item(141, 11).Value = 1714.32
$ws.Cells.Item(141, 13).Value = 3465.68

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4347.836
$ws.Cells.Item(32, 9).Value = 3954.7966
$ws.Cells.Item(32, 11).Value = 3954.7966
$ws.Cells.Item(32, 13).Value = -3667.7966
$ws.Cells.Item(131, 8).Value = 49470
$ws.Cells.Item(131, 10).Value = 49470
$ws.Cells.Item(131, 12).Value = 49470
$ws.Cells.Item(131, 14).Value = -59550

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 76924430
$ws.Cells.Item(16, 9).Value = 100001360
$ws.Cells.Item(16, 10).Value = 1333.3334
$ws.Cells.Item(16, 11).Value = 100001360
$ws.Cells.Item(16, 12).Value = 1333.3334
$ws.Cells.Item(16, 13).Value = -100001073
$ws.Cells.Item(16, 14).Value = -1907.3334
$ws.Cells.Item(31, 8).Value = 1982.6428
$ws.Cells.Item(31, 9).Value = 2100
$ws.Cells.Item(31, 11).Value = 2100
$ws.Cells.Item(31, 13).Value = -1805
$ws.Cells.Item(34, 8).Value = 1982.6428
$ws.Cells.Item(34, 9).Value = 2100
$ws.Cells.Item(34, 11).Value = 2100
$ws.Cells.Item(34, 13).Value = -1898
$ws.Cells.Item(58, 8).Value = 972.1667
$ws.Cells.Item(58, 9).Value = 849.5
$ws.Cells.Item(58, 10).Value = 1340.1666
$ws.Cells.Item(58, 11).Value = 849.5
$ws.Cells.Item(58, 12).Value = 1340.1666
$ws.Cells.Item(58, 13).Value = -646.5
$ws.Cells.Item(58, 14).Value = -1746.1666
$ws.Cells.Item(93, 8).Value = 41666.668
$ws.Cells.Item(93, 10).Value = 41666.668
$ws.Cells.Item(93, 12).Value = 41666.668
$ws.Cells.Item(93, 14).Value = -45410.668
$ws.Cells.Item(99, 8).Value = 1907.7778
$ws.Cells.Item(99, 9).Value = 2203.3333
$ws.Cells.Item(99, 11).Value = 2203.3333
$ws.Cells.Item(99, 13).Value = -705.3332999999998
$ws.Cells.Item(113, 8).Value = 76924430
$ws.Cells.Item(113, 9).Value = 100001360
$ws.Cells.Item(113, 10).Value = 1333.3334
$ws.Cells.Item(113, 11).Value = 100001360
$ws.Cells.Item(113, 12).Value = 1333.3334
$ws.Cells.Item(113, 13).Value = -99999190
$ws.Cells.Item(113, 14).Value = -5673.3334
$ws.Cells.Item(126, 8).Value = 1907.7778
$ws.Cells.Item(126, 9).Value = 2203.3333
$ws.Cells.Item(126, 11).Value = 6609.999899999999
$ws.Cells.Item(126, 13).Value = -4139.999899999999
$ws.Cells.Item(134, 8).Value = 17242428
$ws.Cells.Item(134, 9).Value = 948.381
$ws.Cells.Item(134, 11).Value = 2845.143
$ws.Cells.Item(134, 13).Value = -310.143
$ws.Cells.Item(136, 8).Value = 972.1667
$ws.Cells.Item(136, 9).Value = 849.5
$ws.Cells.Item(136, 10).Value = 1340.1666
$ws.Cells.Item(136, 11).Value = 2548.5
$ws.Cells.Item(136, 12).Value = 4020.4998
$ws.Cells.Item(136, 13).Value = 1.5
$ws.Cells.Item(136, 14).Value = -9120.4998

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 2658.1738
$ws.Cells.Item(39, 10).Value = 2565.158
$ws.Cells.Item(39, 12).Value = 7695.474
$ws.Cells.Item(39, 14).Value = -8283.474
$ws.Cells.Item(55, 8).Value = 3466.6667
$ws.Cells.Item(55, 10).Value = 3466.6667
$ws.Cells.Item(55, 12).Value = 10400.0001
$ws.Cells.Item(55, 14).Value = -10754.0001
$ws.Cells.Item(137, 8).Value = 20275490
$ws.Cells.Item(137, 10).Value = 8604.053
$ws.Cells.Item(137, 12).Value = 25812.159
$ws.Cells.Item(137, 14).Value = -36012.159

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2788.8462
$ws.Cells.Item(80, 9).Value = 1800
$ws.Cells.Item(80, 10).Value = 3228.3333
$ws.Cells.Item(80, 11).Value = 1800
$ws.Cells.Item(80, 12).Value = 3228.3333
$ws.Cells.Item(80, 13).Value = -802
$ws.Cells.Item(80, 14).Value = -5224.3333
$ws.Cells.Item(83, 8).Value = 2788.8462
$ws.Cells.Item(83, 9).Value = 1800
$ws.Cells.Item(83, 10).Value = 3228.3333
$ws.Cells.Item(83, 11).Value = 9000
$ws.Cells.Item(83, 12).Value = 16141.6665
$ws.Cells.Item(83, 13).Value = -4008
$ws.Cells.Item(83, 14).Value = -26125.6665
$ws.Cells.Item(102, 8).Value = 1686.0769
$ws.Cells.Item(102, 9).Value = 1477.375
$ws.Cells.Item(102, 10).Value = 2020
$ws.Cells.Item(102, 11).Value = 1477.375
$ws.Cells.Item(102, 12).Value = 2020
$ws.Cells.Item(102, 13).Value = 144.625
$ws.Cells.Item(102, 14).Value = -5264
$ws.Cells.Item(122, 8).Value = 3056.182
$ws.Cells.Item(122, 9).Value = 2436.5
$ws.Cells.Item(122, 11).Value = 7309.5
$ws.Cells.Item(122, 13).Value = -4859.5
$ws.Cells.Item(126, 8).Value = 2170.2856
$ws.Cells.Item(126, 9).Value = 1837
$ws.Cells.Item(126, 10).Value = 2614.6667
$ws.Cells.Item(126, 11).Value = 5511
$ws.Cells.Item(126, 12).Value = 7844.000100000001
$ws.Cells.Item(126, 13).Value = -3041
$ws.Cells.Item(126, 14).Value = -12784.0001

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 5098
$ws.Cells.Item(46, 9).Value = 1980
$ws.Cells.Item(46, 10).Value = 5444.4443
$ws.Cells.Item(46, 11).Value = 1980
$ws.Cells.Item(46, 12).Value = 5444.4443
$ws.Cells.Item(46, 13).Value = -1792
$ws.Cells.Item(46, 14).Value = -5820.4443
$ws.Cells.Item(100, 8).Value = 1908.1111
$ws.Cells.Item(100, 9).Value = 1654.8
$ws.Cells.Item(100, 11).Value = 1654.8
$ws.Cells.Item(100, 13).Value = -1113.8
$ws.Cells.Item(132, 8).Value = 21194.785
$ws.Cells.Item(132, 9).Value = 1228.75
$ws.Cells.Item(132, 10).Value = 54821.79
$ws.Cells.Item(132, 11).Value = 3686.25
$ws.Cells.Item(132, 12).Value = 164465.37
$ws.Cells.Item(132, 13).Value = -1156.25
$ws.Cells.Item(132, 14).Value = -169525.37
$ws.Cells.Item(136, 8).Value = 1076.4762
$ws.Cells.Item(136, 9).Value = 943.56757
$ws.Cells.Item(136, 11).Value = 2830.70271
$ws.Cells.Item(136, 13).Value = -280.70271

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 465.58823
$ws.Cells.Item(107, 9).Value = 427.3846
$ws.Cells.Item(107, 10).Value = 589.75
$ws.Cells.Item(107, 11).Value = 1282.1538
$ws.Cells.Item(107, 12).Value = 1769.25
$ws.Cells.Item(107, 13).Value = 637.8462
$ws.Cells.Item(107, 14).Value = -5609.25
$ws.Cells.Item(109, 8).Value = 31143.8
$ws.Cells.Item(109, 10).Value = 26344.25
$ws.Cells.Item(109, 12).Value = 26344.25
$ws.Cells.Item(109, 14).Value = -29118.25
$ws.Cells.Item(132, 8).Value = 3664.625
$ws.Cells.Item(132, 9).Value = 4730.1763
$ws.Cells.Item(132, 10).Value = 1076.8572
$ws.Cells.Item(132, 11).Value = 14190.5289
$ws.Cells.Item(132, 12).Value = 3230.5716
$ws.Cells.Item(132, 13).Value = -11660.5289
$ws.Cells.Item(132, 14).Value = -8290.571599999999
